# Add the new "EC" route (München - Buchloe - Memmingen - Lindau-Reutin -
# Bregenz - St. Margrethen SG - St. Gallen - Winterthur - Zürich Flughafen -
# Zürich HB) as the first block of data rows, pushing the existing routes
# (IC1, IC5, IC8, IC81, IR13, IR36, IR75) down by 10 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 blank rows right after the header row (row 1), shifting all
# existing data rows (old row 2 onward) down by 10.
$ws.Rows("2:11").Insert()

# The inserted rows inherit the row-1 header formatting (bold, s="1") from
# the Insert; strip that back to the regular/default style used by the rest
# of the data rows (column E keeps its s="1" style, matching every other
# row in the sheet, so leave that range alone).
$ws.Range("A2:D11").ClearFormats()

$data = @(
  @("EC", "München", 48.140502393602702, 11.5584549856616),
  @("EC", "Buchloe", 48.033469670828303, 10.716229362172699),
  @("EC", "Memmingen", 47.985657959875503, 10.186999441204099),
  @("EC", "Lindau-Reutin", 47.552365076156001, 9.7027184736939098),
  @("EC", "Bregenz", 47.503316777784597, 9.7412019105930305),
  @("EC", "St. Margrethen SG", 47.4531199726677, 9.6393359720897003),
  @("EC", "St. Gallen", 47.423416883991401, 9.3691949225297009),
  @("EC", "Winterthur", 47.500313800000001, 8.7239736000000008),
  @("EC", "Zürich Flughafen", 47.451026916503899, 8.5638494491577095),
  @("EC", "Zürich HB", 47.378844392360598, 8.5366312378669402)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = 2 + $i
  $row = $data[$i]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
}

# Match the saved selection state recorded in the workbook.
$ws.Range("G10").Select()
